# Apply the crypto price/volume refresh described in the commit diff.
# Cells whose new text would otherwise be auto-parsed as a number by Excel
# (e.g. "298.35") are temporarily forced to Text format so the literal
# string is preserved, matching the inlineStr cells in the target XML.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range('D2').Value = '39.630.16'
$ws.Range('E2').Value = '  -1.01%  '
$ws.Range('D3').Value = '2.215.59'
$ws.Range('E3').Value = '  -5.28%  '
$ws.Range('E4').Value = '  +0.02%  '
Set-TextValue 'D5' '298.35'
$ws.Range('E5').Value = '  -3.53%  '
Set-TextValue 'D6' '83.70'
$ws.Range('E6').Value = '  -1.98%  '
$ws.Range('E7').Value = '  -2.92%  '
$ws.Range('E8').Value = '  -0.03%  '
$ws.Range('E9').Value = '  -3.93%  '
$ws.Range('E10').Value = '  -3.81%  '
Set-TextValue 'D11' '29.62'
$ws.Range('E11').Value = '  -1.27%  '
Set-TextValue 'D12' '46.07'
$ws.Range('E12').Value = '  -12.09%  '
$ws.Range('E13').Value = '  -2.22%  '
$ws.Range('D14').Value = '2.557.80'
$ws.Range('E14').Value = '  -4.77%  '
$ws.Range('E15').Value = '  -2.43%  '
Set-TextValue 'D16' '14.11'
$ws.Range('E16').Value = '  -4.44%  '
$ws.Range('D17').Value = '2.209.60'
$ws.Range('E17').Value = '  -6.56%  '
$ws.Range('E18').Value = '  -5.23%  '
$ws.Range('D19').Value = '39.556.97'
$ws.Range('E19').Value = '  -1.15%  '
$ws.Range('E20').Value = '  -2.87%  '
$ws.Range('E21').Value = '  -6.15%  '
Set-TextValue 'D22' '65.00'
$ws.Range('E22').Value = '  -4.33%  '
$ws.Range('E23').Value = '  -2.43%  '
Set-TextValue 'D24' '232.72'
$ws.Range('E25').Value = '  -0.15%  '
$ws.Range('E26').Value = '  -4.90%  '
$ws.Range('E27').Value = '  +0.62%  '
Set-TextValue 'D28' '22.72'
$ws.Range('E28').Value = '  -2.63%  '
Set-TextValue 'D29' '2.18'
$ws.Range('E29').Value = '  +2.51%  '
$ws.Range('E30').Value = '  -1.57%  '
Set-TextValue 'D31' '32.26'
$ws.Range('E31').Value = '  -7.51%  '
Set-TextValue 'D32' '149.37'
$ws.Range('E32').Value = '  -2.57%  '
$ws.Range('E33').Value = '  -0.17%  '
$ws.Range('E34').Value = '  -5.41%  '
Set-TextValue 'D35' '2.38'
$ws.Range('E35').Value = '  -2.37%  '
Set-TextValue 'D36' '0.0701'
$ws.Range('E36').Value = '  -2.58%  '
Set-TextValue 'D37' '16.18'
$ws.Range('E37').Value = '  +3.58%  '
$ws.Range('E38').Value = '  -2.73%  '
Set-TextValue 'D39' '0.0971'
$ws.Range('E39').Value = '  -2.05%  '
$ws.Range('E40').Value = '  -5.49%  '
$ws.Range('E41').Value = '  -4.27%  '
$ws.Range('E42').Value = '  -5.38%  '
$ws.Range('D43').Value = '1.925.94'
$ws.Range('E43').Value = '  -1.56%  '
$ws.Range('E44').Value = '  -3.13%  '
Set-TextValue 'D45' '0.0265'
$ws.Range('E45').Value = '  +0.49%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue 'D46' '16.53'
$ws.Range('E46').Value = '  -5.76%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D47' '9.22'
$ws.Range('E47').Value = '  -1.87%  '
$ws.Range('E48').Value = '  -3.88%  '
$ws.Range('D49').Value = '2.434.69'
$ws.Range('E49').Value = '  -4.93%  '
Set-TextValue 'D50' '70.83'
$ws.Range('E50').Value = '  +0.17%  '
Set-TextValue 'D51' '88.59'
$ws.Range('E51').Value = '  -4.37%  '
